# Refresh the "cryptos" symbol list (GitHub Actions scheduled run):
# updates Price / Volume(1h) figures for the existing rows and swaps in
# the LEO row (and shifts the rows below it down by one) on the coin
# list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (preserve original inlineStr "Text" cell type) for
# the numeric-looking Price (D) / Volume(1h) (E) columns being updated,
# so Excel does not silently reinterpret them as Number/Percentage.
$deCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51")
foreach ($cell in $deCells) { $ws.Range($cell).NumberFormat = "@" }

# Apply the new values cell by cell, in sheet order.
$ws.Range("D2").Value = "325.84"
$ws.Range("E2").Value = "-2.27%"
$ws.Range("D3").Value = "44.55"
$ws.Range("E3").Value = "1.86%"
$ws.Range("D4").Value = "5.597"
$ws.Range("E4").Value = "-3.82%"
$ws.Range("D5").Value = "0.08076"
$ws.Range("E5").Value = "-3.10%"
$ws.Range("D6").Value = "8.679"
$ws.Range("E6").Value = "-1.35%"
$ws.Range("D7").Value = "1.901"
$ws.Range("E7").Value = "-4.37%"
$ws.Range("D8").Value = "4.289"
$ws.Range("E8").Value = "-4.78%"
$ws.Range("D9").Value = "2.670"
$ws.Range("E9").Value = "-7.86%"
$ws.Range("D10").Value = "0.9426"
$ws.Range("E10").Value = "-0.01%"
$ws.Range("D11").Value = "0.1171"
$ws.Range("E11").Value = "-5.75%"
$ws.Range("D12").Value = "0.1864"
$ws.Range("E12").Value = "-4.37%"
$ws.Range("D13").Value = "0.09954"
$ws.Range("E13").Value = "2.47%"
$ws.Range("D14").Value = "0.04272"
$ws.Range("E14").Value = "-7.30%"
$ws.Range("D15").Value = "0.1065"
$ws.Range("E15").Value = "-0.24%"
$ws.Range("D16").Value = "0.001288"
$ws.Range("E16").Value = "-0.71%"
$ws.Range("D17").Value = "0.04196"
$ws.Range("E17").Value = "-4.81%"
$ws.Range("D18").Value = "0.005860"
$ws.Range("E18").Value = "-1.46%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "3.583"
$ws.Range("E19").Value = "2.54%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3501"
$ws.Range("E20").Value = "-0.26%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "8.427"
$ws.Range("E21").Value = "-4.60%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "0.1371"
$ws.Range("E22").Value = "0.62%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2524"
$ws.Range("E23").Value = "-4.15%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "0.001238"
$ws.Range("E24").Value = "-1.53%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "0.004535"
$ws.Range("E25").Value = "2.81%"
$ws.Range("D26").Value = "0.0001182"
$ws.Range("E26").Value = "-1.64%"
$ws.Range("E27").Value = "-0.08%"
$ws.Range("D39").Value = "0.02641"
$ws.Range("E39").Value = "-5.81%"
$ws.Range("D40").Value = "0.05442"
$ws.Range("E40").Value = "-5.33%"
$ws.Range("D41").Value = "0.007703"
$ws.Range("E41").Value = "-2.89%"
$ws.Range("D42").Value = "0.1396"
$ws.Range("E42").Value = "-2.16%"
$ws.Range("D43").Value = "0.007174"
$ws.Range("E43").Value = "-20.67%"
$ws.Range("D44").Value = "0.002027"
$ws.Range("E44").Value = "-4.08%"
$ws.Range("D45").Value = "0.008833"
$ws.Range("E45").Value = "-15.94%"
$ws.Range("D46").Value = "0.00007121"
$ws.Range("E46").Value = "-1.33%"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").Value = "0.003532"
$ws.Range("E48").Value = "8.98%"
$ws.Range("E49").Value = "-0.37%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.01%"
